$d = $word.ActiveDocument

$d.Content.Find.Execute("33×26=858", $true, $false, $false, $false, $false, $true, 1, $false, "41×22=902", 2) | Out-Null
$d.Content.Find.Execute("93×65=6045", $true, $false, $false, $false, $false, $true, 1, $false, "14×50=700", 2) | Out-Null
$d.Content.Find.Execute("59×37=2183", $true, $false, $false, $false, $false, $true, 1, $false, "61×66=4026", 2) | Out-Null
$d.Content.Find.Execute("25×37=925", $true, $false, $false, $false, $false, $true, 1, $false, "44×26=1144", 2) | Out-Null
$d.Content.Find.Execute("89×45=4005", $true, $false, $false, $false, $false, $true, 1, $false, "41×25=1025", 2) | Out-Null
$d.Content.Find.Execute("45×75=3375", $true, $false, $false, $false, $false, $true, 1, $false, "87×52=4524", 2) | Out-Null
$d.Content.Find.Execute("72×18=1296", $true, $false, $false, $false, $false, $true, 1, $false, "15×57=855", 2) | Out-Null
$d.Content.Find.Execute("93×85=7905", $true, $false, $false, $false, $false, $true, 1, $false, "28×13=364", 2) | Out-Null
$d.Content.Find.Execute("73×25=1825", $true, $false, $false, $false, $false, $true, 1, $false, "36×28=1008", 2) | Out-Null
$d.Content.Find.Execute("24×73=1752", $true, $false, $false, $false, $false, $true, 1, $false, "94×91=8554", 2) | Out-Null
$d.Content.Find.Execute("50×64=3200", $true, $false, $false, $false, $false, $true, 1, $false, "21×89=1869", 2) | Out-Null
$d.Content.Find.Execute("93×55=5115", $true, $false, $false, $false, $false, $true, 1, $false, "64×85=5440", 2) | Out-Null
$d.Content.Find.Execute("90×88=7920", $true, $false, $false, $false, $false, $true, 1, $false, "16×62=992", 2) | Out-Null
$d.Content.Find.Execute("73×95=6935", $true, $false, $false, $false, $false, $true, 1, $false, "73×53=3869", 2) | Out-Null
$d.Content.Find.Execute("32×86=2752", $true, $false, $false, $false, $false, $true, 1, $false, "74×14=1036", 2) | Out-Null
$d.Content.Find.Execute("42×72=3024", $true, $false, $false, $false, $false, $true, 1, $false, "94×85=7990", 2) | Out-Null
$d.Content.Find.Execute("20×82=1640", $true, $false, $false, $false, $false, $true, 1, $false, "35×19=665", 2) | Out-Null
$d.Content.Find.Execute("72×84=6048", $true, $false, $false, $false, $false, $true, 1, $false, "45×55=2475", 2) | Out-Null
$d.Content.Find.Execute("22×86=1892", $true, $false, $false, $false, $false, $true, 1, $false, "19×33=627", 2) | Out-Null
$d.Content.Find.Execute("60×57=3420", $true, $false, $false, $false, $false, $true, 1, $false, "15×55=825", 2) | Out-Null
$d.Content.Find.Execute("53×23=1219", $true, $false, $false, $false, $false, $true, 1, $false, "64×75=4800", 2) | Out-Null
$d.Content.Find.Execute("40×14=560", $true, $false, $false, $false, $false, $true, 1, $false, "78×35=2730", 2) | Out-Null
$d.Content.Find.Execute("81×48=3888", $true, $false, $false, $false, $false, $true, 1, $false, "68×84=5712", 2) | Out-Null
$d.Content.Find.Execute("31×65=2015", $true, $false, $false, $false, $false, $true, 1, $false, "71×98=6958", 2) | Out-Null
$d.Content.Find.Execute("31×32=992", $true, $false, $false, $false, $false, $true, 1, $false, "53×14=742", 2) | Out-Null

Write-Host "Replacements complete"
